$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row 6 (copies formatting down from row 5, as Excel's row-insert
# normally does) and fill in the new appropriation entry:
#   Data = 10/09/2013 (serial 41527), Quantidade de horas = 0:40 (40 minutes)
$ws.Rows("6:6").Insert()
$ws.Range("A6").Value = 41527
$ws.Range("B6").Value = 0.027777777777777776

# Insert new (still empty) row 7 below it, ready for the next entry; it only
# keeps the time-formatted cell in column B, so clear out the column A cell
# that the row-insert carried down.
$ws.Rows("7:7").Insert()
$ws.Range("A7").Clear()

# Leave the selection on the new empty cell, as in the edited workbook.
$ws.Range("B7").Select()
